$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 208
$ws.Range("A208").Value = 683
$ws.Range("B208").Value = '683-001'
$ws.Range("C208").Value = 59
$ws.Range("D208").Value = 'Female'
$ws.Range("E208").Value = 'Y'
$ws.Range("F208").Value = 'N'
$ws.Range("G208").Value = 'UK'
$ws.Range("H208").Value = 'UK'
$ws.Range("I208").Value = 'UK'
$ws.Range("J208").Value = 'Y'
$ws.Range("K208").Value = 'N'
$ws.Range("L208").Value = 'Y'
$ws.Range("M208").Value = 'UK'
$ws.Range("N208").Value = 'UK'
$ws.Range("O208").Value = 'UK'
$ws.Range("P208").Value = 'UK'
$ws.Range("Q208").Value = 'Y'
$ws.Range("R208").Value = 'N'
$ws.Range("S208").Value = '"Schizoaffective disorder, hepatitis C, previous substance use disorder, hypertension, and type 2 diabetes"'
$ws.Range("T208").Value = 'N'
$ws.Range("U208").Value = 'N'
$ws.Range("V208").Value = 'Y'
$ws.Range("W208").Value = 'Y'
$ws.Range("X208").Value = 'Y'
$ws.Range("Y208").Value = 'N'
$ws.Range("Z208").Value = '"12cm pen"'
$ws.Range("AA208").Value = 'N'
$ws.Range("AB208").Value = 'Y'
$ws.Range("AC208").Value = 'N'
$ws.Range("AD208").Value = 'Y'
$ws.Range("AE208").Value = 'Y'
$ws.Range("AF208").Value = 'N'
$ws.Range("AG208").Value = 'Y'
$ws.Range("AH208").Value = '"hepatic abscess caused by foreign body migration to present initially with a pericardial effusion."'

# Row 209
$ws.Range("A209").Value = 684
$ws.Range("B209").Value = '684-001'
$ws.Range("C209").Value = 12
$ws.Range("D209").Value = 'Male'
$ws.Range("E209").Value = 'Y'
$ws.Range("F209").Value = 'N'
$ws.Range("G209").Value = 'N'
$ws.Range("H209").Value = 'UK'
$ws.Range("I209").Value = 'N'
$ws.Range("J209").Value = 'UK'
$ws.Range("K209").Value = 'N'
$ws.Range("L209").Value = 'N'
$ws.Range("M209").Value = 'UK'
$ws.Range("N209").Value = 'UK'
$ws.Range("O209").Value = 'UK'
$ws.Range("P209").Value = 'UK'
$ws.Range("Q209").Value = 'Y'
$ws.Range("R209").Value = 'N'
$ws.Range("S209").Value = '"biting and chewing coconut fibre based and plastic fibre based doormat at home for several months and hence a diagnosis of gastric bezoars due to doormat ingestion was made and patient was taken up for surgery after deworming"'
$ws.Range("T209").Value = 'N'
$ws.Range("U209").Value = 'N'
$ws.Range("V209").Value = 'N'
$ws.Range("W209").Value = 'N'
$ws.Range("X209").Value = 'N'
$ws.Range("Y209").Value = 'Y'
$ws.Range("Z209").Value = '"doormat"'
$ws.Range("AA209").Value = 'N'
$ws.Range("AB209").Value = 'Y'
$ws.Range("AC209").Value = 'N'
$ws.Range("AD209").Value = 'Y'
$ws.Range("AE209").Value = 'N'
$ws.Range("AF209").Value = 'N'
$ws.Range("AG209").Value = 'N'

# Row 210
$ws.Range("A210").Value = 686
$ws.Range("B210").Value = '686-001'
$ws.Range("C210").Value = 39
$ws.Range("D210").Value = 'Male'
$ws.Range("E210").Value = 'Y'
$ws.Range("F210").Value = 'N'
$ws.Range("G210").Value = 'N'
$ws.Range("H210").Value = 'UK'
$ws.Range("I210").Value = 'UK'
$ws.Range("J210").Value = 'Y'
$ws.Range("K210").Value = 'N'
$ws.Range("L210").Value = 'Y'
$ws.Range("M210").Value = 'N'
$ws.Range("N210").Value = 'N'
$ws.Range("O210").Value = 'Y'
$ws.Range("P210").Value = 'N'
$ws.Range("Q210").Value = 'N'
$ws.Range("R210").Value = 'N'
$ws.Range("S210").Value = '"When questioned about the reason for swallowing a foreign object, the patient was unable to recall doing so. This unconscious eating behavior, thought to be a direct manifestation of schizophrenia when the patient is stimulated by their external environment."'
$ws.Range("T210").Value = 'N'
$ws.Range("U210").Value = 'N'
$ws.Range("V210").Value = 'Y'
$ws.Range("W210").Value = 'Y'
$ws.Range("X210").Value = 'Y'
$ws.Range("Y210").Value = 'Y'
$ws.Range("Z210").Value = '"120 foreign objects, such as keys, nails, iron bars, needles, nail clippers, blades, and ear spoons, were successfully removed from the patient’s stomach"'
$ws.Range("AA210").Value = 'N'
$ws.Range("AB210").Value = 'Y'
$ws.Range("AC210").Value = 'N'
$ws.Range("AD210").Value = 'Y'
$ws.Range("AE210").Value = 'N'
$ws.Range("AF210").Value = 'N'
$ws.Range("AG210").Value = 'N'

# Row 211
$ws.Range("A211").Value = 692
$ws.Range("B211").Value = '692-001'
$ws.Range("C211").Value = 30
$ws.Range("D211").Value = 'Male'
$ws.Range("E211").Value = 'Y'
$ws.Range("F211").Value = 'Y'
$ws.Range("G211").Value = 'N'
$ws.Range("H211").Value = 'UK'
$ws.Range("I211").Value = 'UK'
$ws.Range("J211").Value = 'N'
$ws.Range("K211").Value = 'N'
$ws.Range("L211").Value = 'N'
$ws.Range("M211").Value = 'N'
$ws.Range("N211").Value = 'N'
$ws.Range("O211").Value = 'N'
$ws.Range("P211").Value = 'Y'
$ws.Range("Q211").Value = 'N'
$ws.Range("R211").Value = 'N'
$ws.Range("S211").Value = '"reason for phone ingestion was to avoid detection and losing the phone to the prison authorities while being detained in prison"'
$ws.Range("T211").Value = 'N'
$ws.Range("U211").Value = 'N'
$ws.Range("V211").Value = 'Y'
$ws.Range("W211").Value = 'Y'
$ws.Range("X211").Value = 'N'
$ws.Range("Y211").Value = 'N'
$ws.Range("Z211").Value = '"cell phone (with the battery in-situ) in two plastic bags before swallowing.", "71.8 mm x 23.5 mm x13.0 mm and weighed about 20 grams"'
$ws.Range("AA211").Value = 'Y'
$ws.Range("AB211").Value = 'N'
$ws.Range("AC211").Value = 'N'
$ws.Range("AD211").Value = 'N'
$ws.Range("AE211").Value = 'N'
$ws.Range("AF211").Value = 'N'
$ws.Range("AG211").Value = 'N'

# Update selection to reflect the next empty row after the new data (mirrors
# the author continuing to the next blank row ready for further entry).
$ws.Range("A212:XFD212").Select()
